$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.905.67"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "3.506.15"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.64%  "

$ws.Range("E7").Value = "  +1.83%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  -1.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.653"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("E12").Value = "  -1.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").Value = "4.059.45"
$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "606.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.51%  "

$ws.Range("D16").Value = "70.032.45"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").Value = "3.500.98"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.992"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "104.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.10%  "

$ws.Range("E25").Value = "  -2.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.83%  "

$ws.Range("E29").Value = "  +5.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +27.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.49%  "

$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "0.0₃0821"
$ws.Range("E35").Value = "  +6.41%  "

$ws.Range("D36").Value = "3.729.38"
$ws.Range("E36").Value = "  +5.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.77%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "498.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.17%  "

$ws.Range("E43").Value = "  +0.42%  "

$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("E47").Value = "  -2.89%  "

$ws.Range("E48").Value = "  +0.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000244"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.57%  "
